$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

# --- 1) Add the new sheet "post1980_Miami_FL.csv" at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "post1980_Miami_FL.csv"

# Match sheetPr / outline settings and page margins used by the rest of the workbook
$ws3.Outline.SummaryRow = 1
$ws3.Outline.SummaryColumn = 1
$ws3.PageSetup.LeftMargin = 0.75 * 72
$ws3.PageSetup.RightMargin = 0.75 * 72
$ws3.PageSetup.TopMargin = 1 * 72
$ws3.PageSetup.BottomMargin = 1 * 72
$ws3.PageSetup.HeaderMargin = 0.5 * 72
$ws3.PageSetup.FooterMargin = 0.5 * 72

# Apply the same header/index style (bold, bordered, centered) used throughout the workbook
$ws2.Range("B1").Copy()
$ws3.Range("B1:J1").PasteSpecial(-4122)
$ws2.Range("A2").Copy()
$ws3.Range("A2:A32").PasteSpecial(-4122)

# --- 2) Populate the post1980_Miami_FL.csv data table ---
# Row 1
$ws3.Cells.Item(1, 2).Value = "Variable Name"
$ws3.Cells.Item(1, 3).Value = "Baseline"
$ws3.Cells.Item(1, 4).Value = ""
$ws3.Cells.Item(1, 5).Value = "ELF"
$ws3.Cells.Item(1, 6).Value = ""
$ws3.Cells.Item(1, 7).Value = "TLF"
$ws3.Cells.Item(1, 8).Value = ""
$ws3.Cells.Item(1, 9).Value = "PP Peak"
$ws3.Cells.Item(1, 10).Value = ""

# Row 2
$ws3.Cells.Item(2, 1).Value = 0
$ws3.Cells.Item(2, 2).Value = "Annual Electrical Demand"
$ws3.Cells.Item(2, 3).Value = 440189.6
$ws3.Cells.Item(2, 4).Value = "kilowatt_hour"
$ws3.Cells.Item(2, 5).Value = "N/A"
$ws3.Cells.Item(2, 6).Value = "N/A"
$ws3.Cells.Item(2, 7).Value = "N/A"
$ws3.Cells.Item(2, 8).Value = "N/A"
$ws3.Cells.Item(2, 9).Value = "N/A"
$ws3.Cells.Item(2, 10).Value = "N/A"

# Row 3
$ws3.Cells.Item(3, 1).Value = 1
$ws3.Cells.Item(3, 2).Value = "Peak Electrical Demand"
$ws3.Cells.Item(3, 3).Value = 112.85
$ws3.Cells.Item(3, 4).Value = "kilowatt"
$ws3.Cells.Item(3, 5).Value = "N/A"
$ws3.Cells.Item(3, 6).Value = "N/A"
$ws3.Cells.Item(3, 7).Value = "N/A"
$ws3.Cells.Item(3, 8).Value = "N/A"
$ws3.Cells.Item(3, 9).Value = "N/A"
$ws3.Cells.Item(3, 10).Value = "N/A"

# Row 4
$ws3.Cells.Item(4, 1).Value = 2
$ws3.Cells.Item(4, 2).Value = "Annual Thermal Demand"
$ws3.Cells.Item(4, 3).Value = 45917.319
$ws3.Cells.Item(4, 4).Value = "kilowatt_hour"
$ws3.Cells.Item(4, 5).Value = "N/A"
$ws3.Cells.Item(4, 6).Value = "N/A"
$ws3.Cells.Item(4, 7).Value = "N/A"
$ws3.Cells.Item(4, 8).Value = "N/A"
$ws3.Cells.Item(4, 9).Value = "N/A"
$ws3.Cells.Item(4, 10).Value = "N/A"

# Row 5
$ws3.Cells.Item(5, 1).Value = 3
$ws3.Cells.Item(5, 2).Value = "Peak Thermal Demand"
$ws3.Cells.Item(5, 3).Value = 90.652
$ws3.Cells.Item(5, 4).Value = "kilowatt"
$ws3.Cells.Item(5, 5).Value = "N/A"
$ws3.Cells.Item(5, 6).Value = "N/A"
$ws3.Cells.Item(5, 7).Value = "N/A"
$ws3.Cells.Item(5, 8).Value = "N/A"
$ws3.Cells.Item(5, 9).Value = "N/A"
$ws3.Cells.Item(5, 10).Value = "N/A"

# Row 6
$ws3.Cells.Item(6, 1).Value = 4
$ws3.Cells.Item(6, 2).Value = "CHP Size"
$ws3.Cells.Item(6, 3).Value = "N/A"
$ws3.Cells.Item(6, 4).Value = "N/A"
$ws3.Cells.Item(6, 5).Value = 40.59
$ws3.Cells.Item(6, 6).Value = "kilowatt"
$ws3.Cells.Item(6, 7).Value = 2.091
$ws3.Cells.Item(6, 8).Value = "kilowatt"
$ws3.Cells.Item(6, 9).Value = 112.85
$ws3.Cells.Item(6, 10).Value = "kilowatt"

# Row 7
$ws3.Cells.Item(7, 1).Value = 5
$ws3.Cells.Item(7, 2).Value = "TES Size"
$ws3.Cells.Item(7, 3).Value = "N/A"
$ws3.Cells.Item(7, 4).Value = "N/A"
$ws3.Cells.Item(7, 5).Value = 25.23
$ws3.Cells.Item(7, 6).Value = "kilowatt_hour"
$ws3.Cells.Item(7, 7).Value = 22.117
$ws3.Cells.Item(7, 8).Value = "kilowatt_hour"
$ws3.Cells.Item(7, 9).Value = 0
$ws3.Cells.Item(7, 10).Value = "kilowatt_hour"

# Row 8
$ws3.Cells.Item(8, 1).Value = 6
$ws3.Cells.Item(8, 2).Value = "Aux Boiler Size"
$ws3.Cells.Item(8, 3).Value = "N/A"
$ws3.Cells.Item(8, 4).Value = "N/A"
$ws3.Cells.Item(8, 5).Value = 90.65000000000001
$ws3.Cells.Item(8, 6).Value = "kilowatt"
$ws3.Cells.Item(8, 7).Value = 90.65000000000001
$ws3.Cells.Item(8, 8).Value = "kilowatt"
$ws3.Cells.Item(8, 9).Value = 90.65000000000001
$ws3.Cells.Item(8, 10).Value = "kilowatt"

# Row 9
$ws3.Cells.Item(9, 1).Value = 7
$ws3.Cells.Item(9, 2).Value = "CHP Electrical Energy Generation"
$ws3.Cells.Item(9, 3).Value = "N/A"
$ws3.Cells.Item(9, 4).Value = "N/A"
$ws3.Cells.Item(9, 5).Value = 317295.05
$ws3.Cells.Item(9, 6).Value = "kilowatt_hour"
$ws3.Cells.Item(9, 7).Value = 15244.17
$ws3.Cells.Item(9, 8).Value = "kilowatt_hour"
$ws3.Cells.Item(9, 9).Value = 719324.34
$ws3.Cells.Item(9, 10).Value = "kilowatt_hour"

# Row 10
$ws3.Cells.Item(10, 1).Value = 8
$ws3.Cells.Item(10, 2).Value = "Electrical Energy Bought"
$ws3.Cells.Item(10, 3).Value = "N/A"
$ws3.Cells.Item(10, 4).Value = "N/A"
$ws3.Cells.Item(10, 5).Value = 122894.55
$ws3.Cells.Item(10, 6).Value = "kilowatt_hour"
$ws3.Cells.Item(10, 7).Value = 424945.44
$ws3.Cells.Item(10, 8).Value = "kilowatt_hour"
$ws3.Cells.Item(10, 9).Value = 61560.76
$ws3.Cells.Item(10, 10).Value = "kilowatt_hour"

# Row 11
$ws3.Cells.Item(11, 1).Value = 9
$ws3.Cells.Item(11, 2).Value = "Electrical Energy Sold"
$ws3.Cells.Item(11, 3).Value = "N/A"
$ws3.Cells.Item(11, 4).Value = "N/A"
$ws3.Cells.Item(11, 5).Value = 0
$ws3.Cells.Item(11, 6).Value = ""
$ws3.Cells.Item(11, 7).Value = 0
$ws3.Cells.Item(11, 8).Value = "kilowatt_hour"
$ws3.Cells.Item(11, 9).Value = 340695.49
$ws3.Cells.Item(11, 10).Value = "kilowatt_hour"

# Row 12
$ws3.Cells.Item(12, 1).Value = 10
$ws3.Cells.Item(12, 2).Value = "CHP Thermal Energy Generation"
$ws3.Cells.Item(12, 3).Value = "N/A"
$ws3.Cells.Item(12, 4).Value = "N/A"
$ws3.Cells.Item(12, 5).Value = 594008.0600000001
$ws3.Cells.Item(12, 6).Value = "kilowatt_hour"
$ws3.Cells.Item(12, 7).Value = 29383.51
$ws3.Cells.Item(12, 8).Value = "kilowatt_hour"
$ws3.Cells.Item(12, 9).Value = 1346647.09
$ws3.Cells.Item(12, 10).Value = "kilowatt_hour"

# Row 13
$ws3.Cells.Item(13, 1).Value = 11
$ws3.Cells.Item(13, 2).Value = "TES Thermal Energy Dispatched"
$ws3.Cells.Item(13, 3).Value = "N/A"
$ws3.Cells.Item(13, 4).Value = "N/A"
$ws3.Cells.Item(13, 5).Value = 83.48999999999999
$ws3.Cells.Item(13, 6).Value = "kilowatt_hour"
$ws3.Cells.Item(13, 7).Value = 2652.33
$ws3.Cells.Item(13, 8).Value = "kilowatt_hour"
$ws3.Cells.Item(13, 9).Value = 0
$ws3.Cells.Item(13, 10).Value = "kilowatt_hour"

# Row 14
$ws3.Cells.Item(14, 1).Value = 12
$ws3.Cells.Item(14, 2).Value = "Boiler Thermal Energy Generation"
$ws3.Cells.Item(14, 3).Value = "N/A"
$ws3.Cells.Item(14, 4).Value = "N/A"
$ws3.Cells.Item(14, 5).Value = 140.13
$ws3.Cells.Item(14, 6).Value = "kilowatt_hour"
$ws3.Cells.Item(14, 7).Value = 17372
$ws3.Cells.Item(14, 8).Value = "kilowatt_hour"
$ws3.Cells.Item(14, 9).Value = 15845.69
$ws3.Cells.Item(14, 10).Value = "kilowatt_hour"

# Row 15
$ws3.Cells.Item(15, 1).Value = 13
$ws3.Cells.Item(15, 2).Value = "CHP Electrical Pct Coverage"
$ws3.Cells.Item(15, 3).Value = "N/A"
$ws3.Cells.Item(15, 4).Value = "N/A"
$ws3.Cells.Item(15, 5).Value = 72.08
$ws3.Cells.Item(15, 6).Value = "%"
$ws3.Cells.Item(15, 7).Value = 3.46
$ws3.Cells.Item(15, 8).Value = "%"
$ws3.Cells.Item(15, 9).Value = 163.41
$ws3.Cells.Item(15, 10).Value = "%"

# Row 16
$ws3.Cells.Item(16, 1).Value = 14
$ws3.Cells.Item(16, 2).Value = "Electricity Bought Pct Coverage"
$ws3.Cells.Item(16, 3).Value = "N/A"
$ws3.Cells.Item(16, 4).Value = "N/A"
$ws3.Cells.Item(16, 5).Value = 27.92
$ws3.Cells.Item(16, 6).Value = "%"
$ws3.Cells.Item(16, 7).Value = 96.54000000000001
$ws3.Cells.Item(16, 8).Value = "%"
$ws3.Cells.Item(16, 9).Value = 13.99
$ws3.Cells.Item(16, 10).Value = "%"

# Row 17
$ws3.Cells.Item(17, 1).Value = 15
$ws3.Cells.Item(17, 2).Value = "CHP Thermal Pct Coverage"
$ws3.Cells.Item(17, 3).Value = "N/A"
$ws3.Cells.Item(17, 4).Value = "N/A"
$ws3.Cells.Item(17, 5).Value = 1293.65
$ws3.Cells.Item(17, 6).Value = "%"
$ws3.Cells.Item(17, 7).Value = 63.99
$ws3.Cells.Item(17, 8).Value = "%"
$ws3.Cells.Item(17, 9).Value = 2932.77
$ws3.Cells.Item(17, 10).Value = "%"

# Row 18
$ws3.Cells.Item(18, 1).Value = 16
$ws3.Cells.Item(18, 2).Value = "TES Thermal Pct Coverage"
$ws3.Cells.Item(18, 3).Value = "N/A"
$ws3.Cells.Item(18, 4).Value = "N/A"
$ws3.Cells.Item(18, 5).Value = 0.18
$ws3.Cells.Item(18, 6).Value = "%"
$ws3.Cells.Item(18, 7).Value = 5.78
$ws3.Cells.Item(18, 8).Value = "%"
$ws3.Cells.Item(18, 9).Value = 0
$ws3.Cells.Item(18, 10).Value = "%"

# Row 19
$ws3.Cells.Item(19, 1).Value = 17
$ws3.Cells.Item(19, 2).Value = "Boiler Thermal Pct Coverage"
$ws3.Cells.Item(19, 3).Value = "N/A"
$ws3.Cells.Item(19, 4).Value = "N/A"
$ws3.Cells.Item(19, 5).Value = 0.31
$ws3.Cells.Item(19, 6).Value = "%"
$ws3.Cells.Item(19, 7).Value = 37.83
$ws3.Cells.Item(19, 8).Value = "%"
$ws3.Cells.Item(19, 9).Value = 34.51
$ws3.Cells.Item(19, 10).Value = "%"

# Row 20
$ws3.Cells.Item(20, 1).Value = 18
$ws3.Cells.Item(20, 2).Value = "Thermal Energy Savings"
$ws3.Cells.Item(20, 3).Value = "N/A"
$ws3.Cells.Item(20, 4).Value = "N/A"
$ws3.Cells.Item(20, 5).Value = -1096970.98
$ws3.Cells.Item(20, 6).Value = "kilowatt_hour"
$ws3.Cells.Item(20, 7).Value = -19770.53
$ws3.Cells.Item(20, 8).Value = "kilowatt_hour"
$ws3.Cells.Item(20, 9).Value = -2579024.68
$ws3.Cells.Item(20, 10).Value = "kilowatt_hour"

# Row 21
$ws3.Cells.Item(21, 1).Value = 19
$ws3.Cells.Item(21, 2).Value = "Electrical Energy Savings"
$ws3.Cells.Item(21, 3).Value = "N/A"
$ws3.Cells.Item(21, 4).Value = "N/A"
$ws3.Cells.Item(21, 5).Value = 793237.63
$ws3.Cells.Item(21, 6).Value = "kilowatt_hour"
$ws3.Cells.Item(21, 7).Value = 38110.42
$ws3.Cells.Item(21, 8).Value = "kilowatt_hour"
$ws3.Cells.Item(21, 9).Value = 946572.12
$ws3.Cells.Item(21, 10).Value = "kilowatt_hour"

# Row 22
$ws3.Cells.Item(22, 1).Value = 20
$ws3.Cells.Item(22, 2).Value = "Total Energy Savings"
$ws3.Cells.Item(22, 3).Value = "N/A"
$ws3.Cells.Item(22, 4).Value = "N/A"
$ws3.Cells.Item(22, 5).Value = -303733.36
$ws3.Cells.Item(22, 6).Value = "kilowatt_hour"
$ws3.Cells.Item(22, 7).Value = 18339.88
$ws3.Cells.Item(22, 8).Value = "kilowatt_hour"
$ws3.Cells.Item(22, 9).Value = -1632452.57
$ws3.Cells.Item(22, 10).Value = "kilowatt_hour"

# Row 23
$ws3.Cells.Item(23, 1).Value = 21
$ws3.Cells.Item(23, 2).Value = "Electricity Cost"
$ws3.Cells.Item(23, 3).Value = 38976.36
$ws3.Cells.Item(23, 4).Value = "dimensionless"
$ws3.Cells.Item(23, 5).Value = 13420.44
$ws3.Cells.Item(23, 6).Value = "dimensionless"
$ws3.Cells.Item(23, 7).Value = 37748.61
$ws3.Cells.Item(23, 8).Value = "dimensionless"
$ws3.Cells.Item(23, 9).Value = 7986.43
$ws3.Cells.Item(23, 10).Value = "dimensionless"

# Row 24
$ws3.Cells.Item(24, 1).Value = 22
$ws3.Cells.Item(24, 2).Value = "Fuel Cost"
$ws3.Cells.Item(24, 3).Value = 8237.23
$ws3.Cells.Item(24, 4).Value = "dimensionless"
$ws3.Cells.Item(24, 5).Value = 33565.12
$ws3.Cells.Item(24, 6).Value = "dimensionless"
$ws3.Cells.Item(24, 7).Value = 8693.709999999999
$ws3.Cells.Item(24, 8).Value = "dimensionless"
$ws3.Cells.Item(24, 9).Value = 67784.16
$ws3.Cells.Item(24, 10).Value = "dimensionless"

# Row 25
$ws3.Cells.Item(25, 1).Value = 23
$ws3.Cells.Item(25, 2).Value = "CHP Installed Cost"
$ws3.Cells.Item(25, 3).Value = "N/A"
$ws3.Cells.Item(25, 4).Value = "N/A"
$ws3.Cells.Item(25, 5).Value = 140184.11
$ws3.Cells.Item(25, 6).Value = "dimensionless"
$ws3.Cells.Item(25, 7).Value = 7222.1
$ws3.Cells.Item(25, 8).Value = "dimensionless"
$ws3.Cells.Item(25, 9).Value = 389719.41
$ws3.Cells.Item(25, 10).Value = "dimensionless"

# Row 26
$ws3.Cells.Item(26, 1).Value = 24
$ws3.Cells.Item(26, 2).Value = "CHP O&M Cost"
$ws3.Cells.Item(26, 3).Value = "N/A"
$ws3.Cells.Item(26, 4).Value = "N/A"
$ws3.Cells.Item(26, 5).Value = 9518.85
$ws3.Cells.Item(26, 6).Value = "dimensionless"
$ws3.Cells.Item(26, 7).Value = 457.32
$ws3.Cells.Item(26, 8).Value = "dimensionless"
$ws3.Cells.Item(26, 9).Value = 21579.73
$ws3.Cells.Item(26, 10).Value = "dimensionless"

# Row 27
$ws3.Cells.Item(27, 1).Value = 25
$ws3.Cells.Item(27, 2).Value = "TES Installed Cost"
$ws3.Cells.Item(27, 3).Value = "N/A"
$ws3.Cells.Item(27, 4).Value = "N/A"
$ws3.Cells.Item(27, 5).Value = 528.77
$ws3.Cells.Item(27, 6).Value = "dimensionless"
$ws3.Cells.Item(27, 7).Value = 463.58
$ws3.Cells.Item(27, 8).Value = "dimensionless"
$ws3.Cells.Item(27, 9).Value = 0
$ws3.Cells.Item(27, 10).Value = "dimensionless"

# Row 28
$ws3.Cells.Item(28, 1).Value = 26
$ws3.Cells.Item(28, 2).Value = "TES O&M Cost"
$ws3.Cells.Item(28, 3).Value = "N/A"
$ws3.Cells.Item(28, 4).Value = "N/A"
$ws3.Cells.Item(28, 5).Value = 0
$ws3.Cells.Item(28, 6).Value = "dimensionless"
$ws3.Cells.Item(28, 7).Value = 0
$ws3.Cells.Item(28, 8).Value = "dimensionless"
$ws3.Cells.Item(28, 9).Value = 0
$ws3.Cells.Item(28, 10).Value = "dimensionless"

# Row 29
$ws3.Cells.Item(29, 1).Value = 27
$ws3.Cells.Item(29, 2).Value = "PP Revenue"
$ws3.Cells.Item(29, 3).Value = "N/A"
$ws3.Cells.Item(29, 4).Value = "N/A"
$ws3.Cells.Item(29, 5).Value = 0
$ws3.Cells.Item(29, 6).Value = "dimensionless"
$ws3.Cells.Item(29, 7).Value = 0
$ws3.Cells.Item(29, 8).Value = "dimensionless"
$ws3.Cells.Item(29, 9).Value = 27323.98
$ws3.Cells.Item(29, 10).Value = "dimensionless"

# Row 30
$ws3.Cells.Item(30, 1).Value = 28
$ws3.Cells.Item(30, 2).Value = "Simple Payback [Yrs]"
$ws3.Cells.Item(30, 3).Value = "N/A"
$ws3.Cells.Item(30, 4).Value = "N/A"
$ws3.Cells.Item(30, 5).Value = -15.15
$ws3.Cells.Item(30, 6).Value = "dimensionless"
$ws3.Cells.Item(30, 7).Value = 24.48
$ws3.Cells.Item(30, 8).Value = "dimensionless"
$ws3.Cells.Item(30, 9).Value = -17.08
$ws3.Cells.Item(30, 10).Value = "dimensionless"

# Row 31
$ws3.Cells.Item(31, 1).Value = 29
$ws3.Cells.Item(31, 2).Value = "Simple Payback (37.5% incentive)"
$ws3.Cells.Item(31, 3).Value = "N/A"
$ws3.Cells.Item(31, 4).Value = "N/A"
$ws3.Cells.Item(31, 5).Value = -9.470000000000001
$ws3.Cells.Item(31, 6).Value = "dimensionless"
$ws3.Cells.Item(31, 7).Value = 15.3
$ws3.Cells.Item(31, 8).Value = "dimensionless"
$ws3.Cells.Item(31, 9).Value = -10.68
$ws3.Cells.Item(31, 10).Value = "dimensionless"

# Row 32
$ws3.Cells.Item(32, 1).Value = 30
$ws3.Cells.Item(32, 2).Value = "CO2"
$ws3.Cells.Item(32, 3).Value = 176
$ws3.Cells.Item(32, 4).Value = "metric_ton"
$ws3.Cells.Item(32, 5).Value = 105
$ws3.Cells.Item(32, 6).Value = "metric_ton"
$ws3.Cells.Item(32, 7).Value = 172
$ws3.Cells.Item(32, 8).Value = "metric_ton"
$ws3.Cells.Item(32, 9).Value = 154
$ws3.Cells.Item(32, 10).Value = "metric_ton"

# --- 3) On the pre1980 sheet, clear the two stray empty-unit cells in row 11 ---
$ws2.Range("F11").ClearContents()
$ws2.Range("H11").ClearContents()

# --- 4) Restore the first sheet as the active tab, matching the workbook state ---
$wb.Worksheets.Item(1).Activate()